$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the "as of" date in the confidentiality notice (cell A16)
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-18 for illustrative purposes only and are subject to change."

# Update the weight (D) and percent-change (E) values for rows 2-13
$ws.Range("D2").Value = 0.03077902646933008
$ws.Range("E2").Value = 0.0006877579092159358

$ws.Range("D3").Value = 0.02361538561251882
$ws.Range("E3").Value = -0.0004733727810650734

$ws.Range("D4").Value = 0.05239065716715507
$ws.Range("E4").Value = -0.003009956008335313

$ws.Range("D5").Value = 0.1375348403495279
$ws.Range("E5").Value = 0.001470828566759286

$ws.Range("D6").Value = 0.03298019566133367
$ws.Range("E6").Value = -0.02292650033715438

$ws.Range("D7").Value = 0.1175043077341324
$ws.Range("E7").Value = -0.007326478149100146

$ws.Range("D8").Value = 0.1033451849290135
$ws.Range("E8").Value = -0.01454545454545453

$ws.Range("D9").Value = 0.03014373305403574
$ws.Range("E9").Value = -0.01227364185110669

$ws.Range("D10").Value = 0.1286284833362727
$ws.Range("E10").Value = -0.01312869822485208

$ws.Range("D11").Value = 0.2413437109816924
$ws.Range("E11").Value = -0.006553693474393829

$ws.Range("D12").Value = 0.1017344747049875
$ws.Range("E12").Value = -0.009033778476040699

$ws.Range("E13").Value = -0.007625066998957042
